# Auto update: 2025-12-06 00:21:02
# Updates the daily hedging/insurance analysis sheet with the new day's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date column (A2:A5) --------------------------------------------------
# Format as text first so Excel does not auto-convert the yyyy-mm-dd string
# into a date serial number; the source data keeps this column as plain text.
$ws.Range("A2:A5").NumberFormat = "@"
$ws.Range("A2:A5").Value = "2025-12-06"
$ws.Range("A2:A5").Style = "Normal"

# --- Row 2: UnitedHealth Group Incorporated (UNH) (B2/C2 unchanged) -------
$ws.Range("D2").Value = 332.22
$ws.Range("E2").Value = 57.9
$ws.Range("F2").Value = 0.74
$ws.Range("H2").Value = 60
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 43
$ws.Range("K2").Value = 57.5
$ws.Range("N2").Value = 51.54219175917372

# --- Row 3: American International Group, Inc. (AIG) ---------------------
$ws.Range("B3").Value = "American International Group, I"
$ws.Range("C3").Value = "AIG"
$ws.Range("D3").Value = 76.81999999999999
$ws.Range("E3").Value = 43.5
$ws.Range("F3").Value = 0.86
$ws.Range("H3").Value = 46
$ws.Range("I3").Value = 50
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 53.5
$ws.Range("N3").Value = 51.54219175917372

# --- Row 4: MetLife, Inc. (MET) -------------------------------------------
$ws.Range("B4").Value = "MetLife, Inc."
$ws.Range("C4").Value = "MET"
$ws.Range("D4").Value = 78.29000000000001
$ws.Range("E4").Value = 48.2
$ws.Range("F4").Value = 2.26
$ws.Range("H4").Value = 23
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 33
$ws.Range("K4").Value = 53.5
$ws.Range("N4").Value = 51.54219175917372

# --- Row 5: Prudential Financial, Inc. (PRU) (B5/C5 unchanged) ------------
$ws.Range("D5").Value = 111.05
$ws.Range("E5").Value = 68.8
$ws.Range("F5").Value = 2.59
$ws.Range("H5").Value = 46
$ws.Range("I5").Value = 40
$ws.Range("J5").Value = 43
$ws.Range("K5").Value = 49.5
$ws.Range("N5").Value = 51.54219175917372
